$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 30: 2025-12-15 (serial 46006), 四方坪站 ---
# Seed the new row by copying formats from the row above (row 29) so the
# existing cell styles (date / currency / integer number formats) are
# reused rather than re-created.
$ws.Range("A29:F29").Copy() | Out-Null
$ws.Range("A30:F30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30, 1).Value = 46006
$ws.Cells.Item(30, 2).Value = "四方坪站"
$ws.Cells.Item(30, 3).Value = 8859.7999999999993
$ws.Cells.Item(30, 4).Value = 7629.86
$ws.Cells.Item(30, 5).Value = 2939.86
$ws.Cells.Item(30, 6).Value = 388

# --- Row 31: 2025-12-15 (serial 46006), 高岭站 ---
$ws.Range("A29:F29").Copy() | Out-Null
$ws.Range("A31:F31").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(31, 1).Value = 46006
$ws.Cells.Item(31, 2).Value = "高岭站"
$ws.Cells.Item(31, 3).Value = 5669.34
$ws.Cells.Item(31, 4).Value = 4697.16
$ws.Cells.Item(31, 5).Value = 1538.96
$ws.Cells.Item(31, 6).Value = 209

$excel.CutCopyMode = $false

# Reflect the scrolled-down view / new active selection from the edit.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("H34").Select() | Out-Null
